$wb = $excel.ActiveWorkbook

# The "survey" sheet is the xlsform definition (first sheet in the workbook).
$survey = $wb.Worksheets.Item(1)

# 1) Clear the stray "mandatory" hint text that used to sit in the hint
#    column (G) of the "type" field (row 4) - it is no longer used.
$survey.Range("G4").Value = ""

# 2) Insert a new hidden "role" field (with default "nurse") right after
#    the "type" field and before the "Nurse Name" field.
$survey.Rows.Item(5).Insert()
$survey.Rows.Item(5).RowHeight = 12.75
$survey.Range("A5").Value = "hidden"
$survey.Range("B5").Value = "role"
$survey.Range("C5").Value = "Role"
$survey.Range("H5").Value = "nurse"

# 3) Give the "nurse" group (row 2) a label of "Nurse".
$survey.Range("C2").Value = "Nurse"

$wb.Save()
